$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the style already used by the
# other header cells (copy format from G1, the "sum" header).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data cell for row 2 (H2) - plain numeric value, no special style.
$ws.Range("H2").Value = 0
